# Applies the "Add files via upload" edit to the 八位序列号收集（收集结果）
# workbook: two earlier submissions get moderated away ("已删除"), and three
# brand-new submissions are appended at the bottom of the response log
# (sheet "八位序列号收集收集结果yd5" / sheet1.xml).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper: force a numeric-looking string into the cell as genuine text
# (t="s" with no quotePrefix styling) by round-tripping it through a text
# formula and then collapsing the formula down to its cached value.
function Set-TextValue($cell, [string]$val) {
    $cell.Formula = '=""&"' + $val + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

# --- 1. Moderate row 120 (伪装 🅥 / f8b9590f / 85229794) -------------------
# Copy the formatting of an existing "deleted" row (row 11) so the
# strikethrough font + grey fill styles match exactly, then blank out the
# QQ-number columns with the placeholder text.
$ws.Range("A11:D11").Copy()
$ws.Range("A120:D120").PasteSpecial(-4122)
$ws.Range("C120").Value = "已删除"
$ws.Range("D120").Value = "已删除"

# --- 2. Moderate row 129 (木 / c135e5d8 / 2358733176) ----------------------
$ws.Range("A11:D11").Copy()
$ws.Range("A129:D129").PasteSpecial(-4122)
$ws.Range("C129").Value = "已删除"
$ws.Range("D129").Value = "已删除"

# --- 3. Append the three new submissions -----------------------------------
$ws.Cells.Item(132, 1).Value = "一念天堂"
$ws.Cells.Item(132, 2).Value = 46005.4341203704
$ws.Cells.Item(132, 2).NumberFormat = "yyyy/m/d h:mm:ss;@"
$ws.Cells.Item(132, 3).Value = "c4ecc993"
Set-TextValue $ws.Cells.Item(132, 4) "906688527"

$ws.Cells.Item(133, 1).Value = "木"
$ws.Cells.Item(133, 2).Value = 46007.244837963
$ws.Cells.Item(133, 2).NumberFormat = "yyyy/m/d h:mm:ss;@"
$ws.Cells.Item(133, 3).Value = "fb8e559 "
Set-TextValue $ws.Cells.Item(133, 4) "2358733476"

$ws.Cells.Item(134, 1).Value = "🦊🐼😺🐯🦁"
$ws.Cells.Item(134, 2).Value = 46010.5293981482
$ws.Cells.Item(134, 2).NumberFormat = "yyyy/m/d h:mm:ss;@"
$ws.Cells.Item(134, 3).Value = "51B51EB5"
Set-TextValue $ws.Cells.Item(134, 4) "851865221"
